# ikea_products.xlsx update
# - Summary sheet: replace PLATSA rows with FANBYN/VIDGA, rename PAX/GRIMO-VIKEDAL,
#   add PAX/TYSSEDAL + NORRFLY rows, move the Total row down accordingly.
# - Czech Data / Poland Data sheets: same product refresh plus two new rows.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth (in characters) is stored in the OOXML <col width> as
# ColumnWidth + 5/6 by this runtime. Helper to get the exact stored width we want.
function Set-ColWidth($ws, $colIndex, $storedWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $storedWidth - 0.8333333333333334
}

# ---------------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Summary")

# Insert two fresh rows right before the old "Total:" row (row 5), pushing it
# down to row 7 while keeping its formatting/formula structure intact.
$ws1.Rows.Item(5).Insert()
$ws1.Rows.Item(5).Insert()

Set-ColWidth $ws1 1 25

# Row 2: PLATSA -> FANBYN / Not available
$ws1.Range("A2").Value = "FANBYN / Not available"
$ws1.Range("B2").Value = "492.284.74"
$ws1.Range("C2").Value = "Not available"
$ws1.Range("D2").Value = 1990
$ws1.Range("E2").Value = "Not available"

# Row 3: PLATSA -> VIDGA
$ws1.Range("A3").Value = "VIDGA"
$ws1.Range("B3").Value = "294.282.52"
$ws1.Range("C3").Value = "Not available"
$ws1.Range("D3").Value = 1246
$ws1.Range("E3").Value = 1030.452245073759

# Row 4: PAX / GRIMO/VIKEDAL -> PAX / GRIMO
$ws1.Range("A4").Value = "PAX / GRIMO"
$ws1.Range("B4").Value = "994.329.72"
$ws1.Range("C4").Value = "200x66x236 cm"
$ws1.Range("D4").Value = 16010
$ws1.Range("E4").Value = 11234.44276946269

# Row 5: new - PAX / TYSSEDAL
$ws1.Range("A5").Value = "PAX / TYSSEDAL"
$ws1.Range("B5").Value = "594.802.72"
$ws1.Range("C5").Value = "150x60x236 cm"
$ws1.Range("D5").Value = 15180
$ws1.Range("E5").Value = 11812.50134596748

# Row 6: new - NORRFLY / Not available
$ws1.Range("A6").Value = "NORRFLY / Not available"
$ws1.Range("B6").Value = "203.322.54"
$ws1.Range("C6").Value = "67 cm"
$ws1.Range("D6").Value = 399
$ws1.Range("E6").Value = "Not available"

# Row 7: Total row, now covering D2:D6 / E2:E6
$ws1.Range("D7").Formula = "=SUM(D2:D6)"
$ws1.Range("E7").Formula = "=SUM(E2:E6)"

# ---------------------------------------------------------------------------
# Sheet 2: Czech Data
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Czech Data")

Set-ColWidth $ws2 1 16
Set-ColWidth $ws2 5 50

# Row 2: PLATSA -> FANBYN
$ws2.Range("A2").Value = "FANBYN"
$ws2.Range("B2").Value = 395.8941347842852
$ws2.Range("C2").Value = 1990
$ws2.Range("D2").Value = "492.284.74"
$ws2.Range("E2").Value = "Židle, bílá/vn./venkovní"
$ws2.Range("F2").Value = "Not available"

# Row 3: PLATSA -> VIDGA
$ws2.Range("A3").Value = "VIDGA"
$ws2.Range("B3").Value = 247.8814532367936
$ws2.Range("C3").Value = 1246
$ws2.Range("D3").Value = "294.282.52"
$ws2.Range("E3").Value = "Sada pro panelové závěsy, nástěnné upevnění/bílá"
$ws2.Range("F3").Value = "Not available"

# Row 4: PAX / GRIMO/VIKEDAL -> PAX / GRIMO
$ws2.Range("A4").Value = "PAX / GRIMO"
$ws2.Range("B4").Value = 3185.057838138897
$ws2.Range("C4").Value = 16010
$ws2.Range("D4").Value = "994.329.72"
$ws2.Range("E4").Value = "Šatní sestava, bílá/bílá,"
$ws2.Range("F4").Value = "200x66x236 cm"

# Row 5: new - PAX / TYSSEDAL
$ws2.Range("A5").Value = "PAX / TYSSEDAL"
$ws2.Range("B5").Value = 3019.936163831884
$ws2.Range("C5").Value = 15180
$ws2.Range("D5").Value = "594.802.72"
$ws2.Range("E5").Value = "Šatní sestava, bílá/zrcadlové sklo,"
$ws2.Range("F5").Value = "150x60x236 cm"

# Row 6: new - NORRFLY
$ws2.Range("A6").Value = "NORRFLY"
$ws2.Range("B6").Value = 79.37776873313055
$ws2.Range("C6").Value = 399
$ws2.Range("D6").Value = "203.322.54"
$ws2.Range("E6").Value = "LED pásek, barva hliníku,"
$ws2.Range("F6").Value = "67 cm"

# ---------------------------------------------------------------------------
# Sheet 3: Poland Data
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Poland Data")

Set-ColWidth $ws3 1 16
Set-ColWidth $ws3 4 15
Set-ColWidth $ws3 5 44

# Row 2: now entirely "Not available"
$ws3.Range("A2").Value = "Not available"
$ws3.Range("B2").Value = "Not available"
$ws3.Range("C2").Value = "Not available"
$ws3.Range("D2").Value = "Not available"
$ws3.Range("E2").Value = "Not available"
$ws3.Range("F2").Value = "Not available"

# Row 3: PLATSA -> VIDGA
$ws3.Range("A3").Value = "VIDGA"
$ws3.Range("B3").Value = 205
$ws3.Range("C3").Value = 1030.452245073759
$ws3.Range("D3").Value = "294.282.52"
$ws3.Range("E3").Value = "Zestaw do zasłon panelowych, ścienna/biały"
$ws3.Range("F3").Value = "Not available"

# Row 4: PAX / GRIMO/VIKEDAL -> PAX / GRIMO
$ws3.Range("A4").Value = "PAX / GRIMO"
$ws3.Range("B4").Value = 2235
$ws3.Range("C4").Value = 11234.44276946269
$ws3.Range("D4").Value = "994.329.72"
$ws3.Range("E4").Value = "Kombinacja szafy, biały/biały,"
$ws3.Range("F4").Value = "200x66x236 cm"

# Row 5: new - PAX / TYSSEDAL
$ws3.Range("A5").Value = "PAX / TYSSEDAL"
$ws3.Range("B5").Value = 2350
$ws3.Range("C5").Value = 11812.50134596748
$ws3.Range("D5").Value = "594.802.72"
$ws3.Range("E5").Value = "Kombinacja szafy, biały/lustro,"
$ws3.Range("F5").Value = "150x60x236 cm"

# Row 6: new - all "Not available"
$ws3.Range("A6").Value = "Not available"
$ws3.Range("B6").Value = "Not available"
$ws3.Range("C6").Value = "Not available"
$ws3.Range("D6").Value = "Not available"
$ws3.Range("E6").Value = "Not available"
$ws3.Range("F6").Value = "Not available"
